# Add sensitivity-test model runs (rows 9-13) and a new Blueprint run (row 14)
# to the all_runs log, matching the style of the existing NGF (row 7) and
# FinalBlueprint (row 8) entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$urbansimPath = '"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.25 - FINAL VERSION"'

$data = @(
  @{Row=9;  A='NGF'; B=2035; C='2035_TM152_NGF_SensDiscount_01'; D='NGF'; E='Sensitivity Test'; F=$urbansimPath; G='run182'; H='current'; I='https://app.asana.com/0/0/1202512897941570/f'},
  @{Row=10; A='NGF'; B=2035; C='2035_TM152_NGF_SensDiscount_02'; D='NGF'; E='Sensitivity Test'; F=$urbansimPath; G='run182'; H='current'; I='https://app.asana.com/0/0/1202512897941573/f'},
  @{Row=11; A='NGF'; B=2035; C='2035_TM152_NGF_SensDiscount_03'; D='NGF'; E='Sensitivity Test'; F=$urbansimPath; G='run182'; H='current'; I='https://app.asana.com/0/0/1202554938897468/f'},
  @{Row=12; A='NGF'; B=2035; C='2035_TM152_NGF_SensDiscount_04'; D='NGF'; E='Sensitivity Test'; F=$urbansimPath; G='run182'; H='current'; I='https://app.asana.com/0/0/1202554938897469/f'},
  @{Row=13; A='NGF'; B=2035; C='2035_TM152_NGF_SensDiscount_05'; D='NGF'; E='Sensitivity Test'; F=$urbansimPath; G='run182'; H='current'; I='https://app.asana.com/0/0/1202554938897470/f'},
  @{Row=14; A='NGF'; B=2035; C='2035_TM152_NGF_Blueprint_00';    D='NGF'; E='Sensitivity Test'; F=$urbansimPath; G='run182'; H='current'; I='https://app.asana.com/0/0/1202521542566668/f'}
)

foreach ($r in $data) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
}

Write-Host "values set"
